$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "AddCustomerTest"

$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"

$ws.Range("A2").Value = "Raman"
$ws.Range("B2").Value = "Arora"
$ws.Range("C2").Value = "A234wd"

$ws.Range("C9").Select()
